$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 694, shifting existing rows 694:735 down to 695:736.
$ws.Rows.Item(694).Insert()

# Column A holds a date-like string ("yyyy/mm/dd") that must stay literal text,
# not get auto-converted to a date serial number - force text format first,
# then reset the style so the cell doesn't end up with a lingering custom format.
$ws.Range("A694").NumberFormat = "@"
$ws.Range("A694").Value = "2026/01/22"
$ws.Range("A694").Style = "Normal"
$ws.Range("B694").Value = "木"
$ws.Range("C694").Value = 19
$ws.Range("D694").Value = 201

Write-Host "done"
